# Case1/case1.xlsx edit: rework "Test Case 4/5" rows into two new scenarios
# ("Passing All Characters" / "Passing All Integers"), merge their
# Description cells, and restyle the header row + Result column with the
# built-in "Good" (green) cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Rewrite the data for rows 14-21 (test cases 4 and 5).
#    Test case 4 becomes "Passing All Characters" (rows 14-17).
#    Test case 5 becomes "Passing All Integers" (rows 18-21).
# ---------------------------------------------------------------------

# Test case 4 header / first data row
$ws.Range("B14").Value = "Passing All Characters"
$ws.Range("C14").Value = "                       A"
$ws.Range("D14").Value = "A Char 1"
$ws.Range("E14").Value = "A Char 1"

# Row 15
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "                       A"
$ws.Range("D15").Value = "A Char 1"
$ws.Range("E15").Value = "A Char 1"

# Row 16
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "                      A"
$ws.Range("D16").Value = "A Char 1"
$ws.Range("E16").Value = "A Char 1"

# Row 17
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "                      A"
$ws.Range("D17").Value = "1 Integer 4"
$ws.Range("E17").Value = "1 Integer 4"

# Test case 5 header / first data row
$ws.Range("B18").Value = "Passing All Integers"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "1 Integer 4"
$ws.Range("E18").Value = "1 Integer 4"

# Row 19
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = "1 Integer 4"
$ws.Range("E19").Value = "1 Integer 4"

# Row 20
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "1 Integer 4"
$ws.Range("E20").Value = "1 Integer 4"

# Row 21
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = "1 Integer 4"
$ws.Range("E21").Value = "1 Integer 4"

# ---------------------------------------------------------------------
# 2. Merge the new "Description" cells for the two regrouped test cases
#    (mirrors the pre-existing A/F column grouping).
# ---------------------------------------------------------------------

$ws.Range("B14:B17").Merge()
$ws.Range("B18:B21").Merge()

# The merged description cells keep the same font as the other
# description cells but are centred horizontally only.
$ws.Range("B14:B17").HorizontalAlignment = -4108
$ws.Range("B18:B21").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3. Apply the built-in "Good" style (green fill / green text) to the
#    header row and to the "Result" column.
# ---------------------------------------------------------------------

$ws.Range("A1:F1").Style = "Good"

$ws.Range("F2:F21").Style = "Good"
$ws.Range("F2:F21").HorizontalAlignment = -4108
$ws.Range("F2:F21").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 4. Misc bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------

$ws.Range("E15").Select()
